# Regenerate the handback status report: refresh the "last generated"
# timestamps that are stamped onto the Overview / per-locale sheets each
# time the report runs.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for bf640afe-...md
$wsOverview.Range("G2").Value = "2016-08-31 05:06:09"

# zh-cn sheet, row for bf640afe-...md:
#   Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-31 05:05:58"
$wsZhCn.Range("K2").Value = "2016-08-31 05:06:47"

# de-de sheet, row for bf640afe-...md:
#   Correspond Handoff Datetime (mirrors the Overview "Latest HO Xliff
#   Generate Date" value) / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-08-31 05:06:09"
$wsDeDe.Range("K2").Value = "2016-08-31 05:06:54"
